$d = $word.ActiveDocument

$pairs = @(
    @("44×17=748", "37×73=2701"),
    @("11×94=1034", "74×51=3774"),
    @("72×48=3456", "30×32=960"),
    @("12×36=432", "78×33=2574"),
    @("86×34=2924", "26×44=1144"),
    @("52×68=3536", "83×24=1992"),
    @("12×66=792", "98×85=8330"),
    @("12×48=576", "89×50=4450"),
    @("82×23=1886", "11×20=220"),
    @("31×96=2976", "80×96=7680"),
    @("56×29=1624", "57×56=3192"),
    @("72×64=4608", "34×30=1020"),
    @("46×59=2714", "73×58=4234"),
    @("80×60=4800", "34×36=1224"),
    @("95×89=8455", "84×28=2352"),
    @("74×80=5920", "12×86=1032"),
    @("65×20=1300", "23×30=690"),
    @("74×65=4810", "68×39=2652"),
    @("24×49=1176", "65×25=1625"),
    @("27×21=567", "99×70=6930"),
    @("78×45=3510", "18×93=1674"),
    @("17×36=612", "55×86=4730"),
    @("81×47=3807", "63×30=1890"),
    @("86×17=1462", "40×87=3480"),
    @("23×51=1173", "12×91=1092")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
